$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (first test user)
$ws.Range("A2").Value = "FName1yy"
$ws.Range("B2").Value = "LName1yy"
$ws.Range("C2").Value = "User#1yy"
$ws.Range("D2").Value = "Admin"
$ws.Range("E2").Value = "admin@mail.com"
$ws.Range("F2").Value = "082555yy"
$ws.Range("G2").Value = "Pass1yy"
$ws.Range("H2").Value = "Company AAA"

# Row 3 (second test user) - same username as row 2, to exercise the
# "unique username" validation being added.
$ws.Range("A3").Value = "FName2yy"
$ws.Range("B3").Value = "LName2yy"
$ws.Range("C3").Value = "User#1yy"
$ws.Range("D3").Value = "Customer"
$ws.Range("E3").Value = "cusomter@mail.com"
$ws.Range("F3").Value = "083444yy"
$ws.Range("G3").Value = "Pass1yy"
$ws.Range("H3").Value = "Company BBB"

$ws.Range("C6").Select()
